# Update cryptos list: refreshed Price (D) and Volume(1h) (E) values.
# Numeric-looking Price values are forced back to text (matching the
# original inlineStr storage) via NumberFormat "@" + Style "Normal" so
# Excel doesn't silently convert them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.311.81"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.65"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.356"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0739"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0975"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("D14").Value = "2.152.27"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.764"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "1.894.31"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "35.352.22"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("E24").Value = "  +7.11%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").Value = "4.128.45"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +5.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.837"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0723"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").Value = "1.303.17"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0798"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("E51").Value = "  -5.74%  "
